$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above current row 51 (shifts existing rows 51-59 down to 53-61)
$ws.Rows("51:52").Insert()

# ---- New row 51 ----
$ws.Cells.Item(51, 1).Value = 10
$ws.Cells.Item(51, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(51, 3).Value = "La Araucanía"
$ws.Cells.Item(51, 4).Value = 44889
$ws.Cells.Item(51, 5).Value = 9
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100103
$ws.Cells.Item(51, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(51, 9).Value = 100103003
$ws.Cells.Item(51, 10).Value = "Damasco"
$ws.Cells.Item(51, 11).Value = "Castle Brite"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 15000
$ws.Cells.Item(51, 15).Value = 15000
$ws.Cells.Item(51, 16).Value = 15000
$ws.Cells.Item(51, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(51, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(51, 19).Value = 1500
$ws.Cells.Item(51, 20).Value = 10

# ---- New row 52 ----
$ws.Cells.Item(52, 1).Value = 10
$ws.Cells.Item(52, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value = "La Araucanía"
$ws.Cells.Item(52, 4).Value = 44889
$ws.Cells.Item(52, 5).Value = 9
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100103
$ws.Cells.Item(52, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(52, 9).Value = 100103003
$ws.Cells.Item(52, 10).Value = "Damasco"
$ws.Cells.Item(52, 11).Value = "Castle Brite"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 90
$ws.Cells.Item(52, 14).Value = 33000
$ws.Cells.Item(52, 15).Value = 34000
$ws.Cells.Item(52, 16).Value = 33556
$ws.Cells.Item(52, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(52, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 19).Value = 1864
$ws.Cells.Item(52, 20).Value = 18

Write-Host "done"
